# Modificar categorias de densidad poblacional
# Collapse the 5-tier Density Class classification (Very High / High / Medium /
# Low / Very Low) into a 3-tier one (High / Medium / Low):
#   Very High -> High
#   High      -> High   (unchanged)
#   Medium    -> High
#   Low       -> Medium
#   Very Low  -> Low

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E2").Value  = "Low"
$ws.Range("E3").Value  = "Medium"
$ws.Range("E4").Value  = "Low"
$ws.Range("E5").Value  = "Medium"
$ws.Range("E7").Value  = "High"
$ws.Range("E8").Value  = "Low"
$ws.Range("E9").Value  = "High"
$ws.Range("E11").Value = "High"
$ws.Range("E12").Value = "Medium"
$ws.Range("E13").Value = "Medium"
$ws.Range("E14").Value = "High"
$ws.Range("E15").Value = "Low"
$ws.Range("E16").Value = "Low"
$ws.Range("E17").Value = "Medium"
$ws.Range("E18").Value = "Medium"
$ws.Range("E19").Value = "High"
$ws.Range("E21").Value = "Medium"
$ws.Range("E22").Value = "Medium"
$ws.Range("E23").Value = "Low"
$ws.Range("E24").Value = "Low"
$ws.Range("E25").Value = "Low"
$ws.Range("E26").Value = "Low"
$ws.Range("E27").Value = "Medium"
$ws.Range("E29").Value = "High"
$ws.Range("E30").Value = "Medium"
$ws.Range("E31").Value = "Low"
$ws.Range("E32").Value = "High"
$ws.Range("E33").Value = "Medium"
$ws.Range("E34").Value = "Medium"
$ws.Range("E35").Value = "High"
$ws.Range("E36").Value = "Low"

# Update the selected cell to match the author's final cursor position.
$ws.Range("H16").Select() | Out-Null
